$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New prediction values for rows 2..156 (row 1 is the "Prediction" header).
$values = @(0.1056026443839073,0.09806323796510696,0.09596473723649979,0.09087328612804413,0.09317360073328018,0.1015491411089897,0.1451462507247925,0.1764229983091354,0.1978478878736496,0.1990372389554977,0.1978534013032913,0.1886269897222519,0.1854834705591202,0.1762103885412216,0.1636193096637726,0.162919819355011,0.152069017291069,0.1397820264101028,0.1369117349386215,0.1410485506057739,0.1395169049501419,0.1311882734298706,0.1133195608854294,0.1062970608472824,0.0944255068898201,0.08999232947826385,0.08739485591650009,0.08625379949808121,0.08714801073074341,0.09049210697412491,0.09036631137132645,0.08960402756929398,0.08972924947738647,0.08797987550497055,0.08824003487825394,0.08589421957731247,0.08397543430328369,0.08113499730825424,0.07932413369417191,0.07689343392848969,0.07447367906570435,0.07485280185937881,0.07302823662757874,0.07566092908382416,0.07618246227502823,0.07336227595806122,0.07110990583896637,0.07004912942647934,0.08768191188573837,0.08452703058719635,0.08288638293743134,0.08208984136581421,0.08298405259847641,0.08703645318746567,0.08487524092197418,0.0821402370929718,0.08285330981016159,0.08223347365856171,0.08016877621412277,0.07626821845769882,0.07129251956939697,0.07426422834396362,0.07452625781297684,0.07611998915672302,0.07241465151309967,0.07692595571279526,0.07759345322847366,0.08485069125890732,0.08583956211805344,0.08014824986457825,0.07356823980808258,0.07434379309415817,0.0879397988319397,0.08594207465648651,0.08859008550643921,0.08848121017217636,0.08726364374160767,0.0948173776268959,0.1319430023431778,0.1777654439210892,0.2061997503042221,0.2222494632005692,0.2104895412921906,0.2053017914295197,0.2013587802648544,0.1938804686069489,0.1848463267087936,0.1858280748128891,0.1699197292327881,0.1678435504436493,0.1636290103197098,0.1648920774459839,0.1667213886976242,0.1615720093250275,0.1247084885835648,0.1073317229747772,0.1046975776553154,0.09925060719251633,0.1045616045594215,0.1037989035248756,0.09332532435655594,0.1003702506422997,0.1444872915744781,0.1839158833026886,0.2095482051372528,0.2254806011915207,0.2193359732627869,0.2065916657447815,0.2056795954704285,0.1968738734722137,0.1876902282238007,0.1876541525125504,0.1775341182947159,0.1671268790960312,0.1628993004560471,0.1719160526990891,0.168288066983223,0.1682065725326538,0.1454377472400665,0.1245240643620491,0.1087479367852211,0.1011644825339317,0.09873849898576736,0.09690376371145248,0.0977586954832077,0.1050778105854988,0.1460853666067123,0.1909648329019547,0.2153187543153763,0.220765620470047,0.220290794968605,0.2097830921411514,0.2113891690969467,0.202499195933342,0.1886700540781021,0.1873032301664352,0.177996352314949,0.1672581732273102,0.1639053523540497,0.1671480089426041,0.1664849519729614,0.1615317910909653,0.1355683654546738,0.115163654088974,0.1091917157173157,0.1016082614660263,0.09950978308916092,0.09830059856176376,0.09899407625198364,0.1082462519407272,0.1490663141012192,0.1841649115085602,0.2147525697946548,0.2201407104730606,0.2224241495132446)

# Capture the existing number format used by the original data column (A2)
# before we start overwriting, so newly-created rows (96..156) match the
# same style as the pre-existing ones (s="1" -> numFmt "0.000").
$dataFormat = $ws.Range("A2").NumberFormat

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $values[$i]
    $cell.NumberFormat = $dataFormat
}
